# AHB-Diff header rework:
#   - "<field>_old" headers (columns A-J) -> "<field>_FV2404"
#   - "<field>_new" headers (columns L-U) -> "<field>_FV2410"
#   - promote the data range A1:U78 to a real Excel Table ("Table1")
#   - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

# Columns 1-10 (A-J) hold the "_old" -> "_FV2404" headers.
for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2404Headers[$i]
}

# Column 11 (K) is the untouched "diff" header.
# Columns 12-21 (L-U) hold the "_new" -> "_FV2410" headers.
for ($i = 0; $i -lt $fv2410Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2410Headers[$i]
}

# Convert the used range into an Excel Table so the header row gets filter
# dropdowns and the workbook carries a proper ListObject named "Table1".
$rng = $ws.Range("A1:U78")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1, top pane frozen).
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
